$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: per-source-file info needed for "Latest Target File" (col I)
# and "Latest Handback File" (col J) on the zh-cn / de-de sheets, plus the
# matching hyperlink target URLs (same URLs already used by column A).
# ---------------------------------------------------------------------------
$ymlUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3041b1f571689ec9c613a460e17daff7f6cc594d/e2e/96590062-658a-458a-8dd4-c9bb9de1dfc7.yml"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3041b1f571689ec9c613a460e17daff7f6cc594d/e2e/c04a5d1a-3095-42da-bb42-89b2045bb8fd.md"
$yml2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3041b1f571689ec9c613a460e17daff7f6cc594d/e2e/f36eb4b8-7a6b-4b0e-ac2a-cb450733d80b.yml"

$ymlName = "96590062-658a-458a-8dd4-c9bb9de1dfc7.yml"
$mdName  = "c04a5d1a-3095-42da-bb42-89b2045bb8fd.md"
$yml2Name = "f36eb4b8-7a6b-4b0e-ac2a-cb450733d80b.yml"

# ---------------------------------------------------------------------------
# 1) Overview sheet: handback status text changed
#    "In Translation" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus

# Widen status columns to fit the longer text (best effort; Excel quantizes
# ColumnWidth to 1/6 character increments, so we pick the closest match).
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: fill in Latest Target File (I) / Latest Handback File (J)
#    and refresh Latest Handback DateTime (K)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Value = $ymlName
$wsZh.Range("J2").Value = "96590062-658a-458a-8dd4-c9bb9de1dfc7.45c6f3e33dc1b13c35e8f8551639bdad3d93a024.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-26 09:00:48"

$wsZh.Range("I3").Value = $mdName
$wsZh.Range("J3").Value = "c04a5d1a-3095-42da-bb42-89b2045bb8fd.9fdf738a85602eb0265dacd025fae9f34da68b27.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-26 09:00:48"

$wsZh.Range("I4").Value = $yml2Name
$wsZh.Range("J4").Value = "f36eb4b8-7a6b-4b0e-ac2a-cb450733d80b.b51cba0b0aa0f502df6fc7aff631961f1ecf4bd7.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-10-26 09:00:48"

# Give the newly populated "Latest Target File" cells the same hyperlink
# look & link as the matching "Source File Name" cells in column A.
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276
$wsZh.Range("I4").Font.Underline = 2
$wsZh.Range("I4").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ymlUrl, [Type]::Missing, [Type]::Missing, $ymlName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $yml2Url, [Type]::Missing, [Type]::Missing, $yml2Name)

# Column width adjustments (Status column C, and new Target/Handback columns I, J)
$wsZh.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 3) de-de sheet: same pattern as zh-cn, but with de-de xlf file names and a
#    distinct handback timestamp
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Value = $ymlName
$wsDe.Range("J2").Value = "96590062-658a-458a-8dd4-c9bb9de1dfc7.45c6f3e33dc1b13c35e8f8551639bdad3d93a024.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-26 09:01:06"

$wsDe.Range("I3").Value = $mdName
$wsDe.Range("J3").Value = "c04a5d1a-3095-42da-bb42-89b2045bb8fd.9fdf738a85602eb0265dacd025fae9f34da68b27.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-26 09:01:06"

$wsDe.Range("I4").Value = $yml2Name
$wsDe.Range("J4").Value = "f36eb4b8-7a6b-4b0e-ac2a-cb450733d80b.b51cba0b0aa0f502df6fc7aff631961f1ecf4bd7.de-de.xlf"
$wsDe.Range("K4").Value = "2016-10-26 09:01:06"

$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276
$wsDe.Range("I4").Font.Underline = 2
$wsDe.Range("I4").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ymlUrl, [Type]::Missing, [Type]::Missing, $ymlName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $yml2Url, [Type]::Missing, [Type]::Missing, $yml2Name)

$wsDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated"
